# Ctenophora.xlsx / "Materials" sheet - add suborder / infraorder / superfamily
# columns (Darwin Core higher-classification ranks) right after "order", and
# populate a couple of previously-blank template cells in row 2.
#
# Before:  ... AQ=order  AR=family      AS=genus  AT=subgenus ...
# After:   ... AQ=order  AR=suborder  AS=infraorder  AT=superfamily  AU=family  AV=genus  AW=subgenus ...

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Insert three new blank columns right before the existing "family" column (AR),
# pushing family/genus/subgenus/... three places to the right.
$ws.Range("AR:AT").EntireColumn.Insert()

# Row 1 holds the Darwin Core field names (the column headers).
$ws.Range("AR1").Value = "suborder"
$ws.Range("AS1").Value = "infraorder"
$ws.Range("AT1").Value = "superfamily"

# Row 2 holds the template/value row with the matching placeholders.
# scientificName's template now pulls from the iNaturalist summary object.
$ws.Range("AG2").Value = "`${summary.taxonName}"

$ws.Range("AR2").Value = "`${suborder}"
$ws.Range("AS2").Value = "`${infraorder}"
$ws.Range("AT2").Value = "`${superfamily}"

# scientificNameAuthorship was blank before - now templated too.
$ws.Range("BB2").Value = "`${summary.Author}"

# eventTime was blank before - now templated with a time-only date format.
$ws.Range("EA2").Value = "!Date:HH:mm:ss"
